# Add a "Responsible by" field to the Q&A sheet:
#   - insert a new header row at the top (TOR # / Description / Responsible by)
#   - populate the new column C with "Excise" / "Huawei" for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing rows down by one to make room for the header row.
$ws.Rows.Item(1).Insert()

# Fill in "Responsible by" (column C) for each data row, in sheet order,
# before writing the header text (matches how the shared-string table ends
# up ordered: Excise, Huawei, TOR #, Description, Responsible by).
$ws.Range("C3").Value = "Excise"
$ws.Range("C4").Value = "Huawei"
$ws.Range("C8").Value = "Excise"
$ws.Range("C9").Value = "Excise"
$ws.Range("C10").Value = "Excise"
$ws.Range("C11").Value = "Excise"
$ws.Range("C12").Value = "Excise"
$ws.Range("C13").Value = "Huawei"
$ws.Range("C14").Value = "Excise"

# New header row.
$ws.Range("A1").Value = "TOR #"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Responsible by"

# Match the author's final selection.
$ws.Range("C15").Select()
